# Update the "取得日時" (acquired datetime) timestamps in the "ランサーズ" sheet
# from 2025-12-12 06:29:31 to 2025-12-12 06:38:41 for all data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-12-12 06:29:31"
$newValue = "2025-12-12 06:38:41"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
